$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 25.88889
$ws.Range("I5").Value = 23.875
$ws.Range("J5").Value = 42
$ws.Range("K5").Value = 23.875
$ws.Range("L5").Value = 42
$ws.Range("M5").Value = 91.125
$ws.Range("N5").Value = -272

$ws.Range("H17").Value = 1209.1111
$ws.Range("J17").Value = 1209.1111
$ws.Range("L17").Value = 3627.3333
$ws.Range("N17").Value = -3963.3333

$ws.Range("H32").Value = 861.3077
$ws.Range("I32").Value = 117
$ws.Range("J32").Value = 1084.6
$ws.Range("K32").Value = 117
$ws.Range("L32").Value = 1084.6
$ws.Range("M32").Value = 209
$ws.Range("N32").Value = -1736.6

$ws.Range("H40").Value = 994.1923
$ws.Range("I40").Value = 969.5
$ws.Range("J40").Value = 996.25
$ws.Range("K40").Value = 969.5
$ws.Range("L40").Value = 996.25
$ws.Range("M40").Value = -794.5
$ws.Range("N40").Value = -1346.25

$ws.Range("H113").Value = 4717.5
$ws.Range("I113").Value = 4561
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 4561
$ws.Range("L113").Value = 5500
$ws.Range("M113").Value = -1307
$ws.Range("N113").Value = -12008

$ws.Range("H137").Value = 1925717
$ws.Range("I137").Value = 4001993.5
$ws.Range("J137").Value = 3238.7407
$ws.Range("K137").Value = 12005980.5
$ws.Range("L137").Value = 9716.222099999999
$ws.Range("M137").Value = -12003430.5
$ws.Range("N137").Value = -14816.2221

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2340.3215
$ws.Range("I2").Value = 2363.5625
$ws.Range("J2").Value = 2309.3333
$ws.Range("K2").Value = 2363.5625
$ws.Range("L2").Value = 2309.3333
$ws.Range("M2").Value = -2250.5625
$ws.Range("N2").Value = -2535.3333

$ws.Range("H32").Value = 25394.146
$ws.Range("I32").Value = 22198.145
$ws.Range("J32").Value = 33875.848
$ws.Range("K32").Value = 22198.145
$ws.Range("L32").Value = 33875.848
$ws.Range("M32").Value = -21911.145
$ws.Range("N32").Value = -34449.848

$ws.Range("H116").Value = 2340.3215
$ws.Range("I116").Value = 2363.5625
$ws.Range("J116").Value = 2309.3333
$ws.Range("K116").Value = 2363.5625
$ws.Range("L116").Value = 2309.3333
$ws.Range("M116").Value = -69.5625
$ws.Range("N116").Value = -6897.3333

$ws.Range("H122").Value = 1932
$ws.Range("I122").Value = 1400.2858
$ws.Range("J122").Value = 2552.3333
$ws.Range("K122").Value = 4200.857400000001
$ws.Range("L122").Value = 7656.999899999999
$ws.Range("M122").Value = -1750.857400000001
$ws.Range("N122").Value = -12556.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2340.3215
$ws.Range("I3").Value = 2363.5625
$ws.Range("J3").Value = 2309.3333
$ws.Range("K3").Value = 2363.5625
$ws.Range("L3").Value = 2309.3333
$ws.Range("M3").Value = -2249.5625
$ws.Range("N3").Value = -2537.3333

$ws.Range("H20").Value = 1771.5385
$ws.Range("I20").Value = 1759.2
$ws.Range("J20").Value = 1779.25
$ws.Range("K20").Value = 1759.2
$ws.Range("L20").Value = 1779.25
$ws.Range("M20").Value = -1512.2
$ws.Range("N20").Value = -2273.25

$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("N57").ClearContents()

$ws.Range("H99").Value = 963.5
$ws.Range("I99").Value = 973.5454999999999
$ws.Range("J99").Value = 926.6667
$ws.Range("K99").Value = 973.5454999999999
$ws.Range("L99").Value = 926.6667
$ws.Range("M99").Value = 524.4545000000001
$ws.Range("N99").Value = -3922.6667

$ws.Range("H105").Value = 2116.1538
$ws.Range("I105").Value = 2042.5
$ws.Range("K105").Value = 2042.5
$ws.Range("M105").Value = -295.5

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1533.6666
$ws.Range("I16").Value = 1050.5
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 1050.5
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = -763.5
$ws.Range("N16").Value = -3074

$ws.Range("H22").Value = 1041947.1
$ws.Range("I22").Value = 1602776
$ws.Range("J22").Value = 407.85715
$ws.Range("K22").Value = 1602776
$ws.Range("L22").Value = 407.85715
$ws.Range("M22").Value = -1602426
$ws.Range("N22").Value = -1107.85715

$ws.Range("H113").Value = 1533.6666
$ws.Range("I113").Value = 1050.5
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1050.5
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 1119.5
$ws.Range("N113").Value = -6840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 3772.9
$ws.Range("J103").Value = 5114.857
$ws.Range("L103").Value = 15344.571
$ws.Range("N103").Value = -17102.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 119411.11
$ws.Range("I70").Value = 504900
$ws.Range("J70").Value = 9271.429
$ws.Range("K70").Value = 504900
$ws.Range("L70").Value = 9271.429
$ws.Range("M70").Value = -504630
$ws.Range("N70").Value = -9811.429

$ws.Range("H73").Value = 119411.11
$ws.Range("I73").Value = 504900
$ws.Range("J73").Value = 9271.429
$ws.Range("K73").Value = 504900
$ws.Range("L73").Value = 9271.429
$ws.Range("M73").Value = -503964
$ws.Range("N73").Value = -11143.429

$ws.Range("H113").Value = 2437.7778
$ws.Range("I113").Value = 991.6667
$ws.Range("J113").Value = 5330
$ws.Range("K113").Value = 991.6667
$ws.Range("L113").Value = 5330
$ws.Range("M113").Value = 1178.3333
$ws.Range("N113").Value = -9670

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2755510
$ws.Range("I46").Value = 4329580
$ws.Range("J46").Value = 887.5
$ws.Range("K46").Value = 4329580
$ws.Range("L46").Value = 887.5
$ws.Range("M46").Value = -4329392
$ws.Range("N46").Value = -1263.5

$ws.Range("H64").Value = 31211.334
$ws.Range("J64").Value = 31211.334
$ws.Range("L64").Value = 31211.334
$ws.Range("N64").Value = -31661.334

$ws.Range("H67").Value = 31211.334
$ws.Range("J67").Value = 31211.334
$ws.Range("L67").Value = 31211.334
$ws.Range("N67").Value = -32771.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 975.5
$ws.Range("I122").Value = 971.38464
$ws.Range("J122").Value = 993.3333
$ws.Range("K122").Value = 2914.15392
$ws.Range("L122").Value = 2979.9999
$ws.Range("M122").Value = -464.1539199999997
$ws.Range("N122").Value = -7879.9999

$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -530
